$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# UPDATED 12/30/2019: P1219 Solved
#   - adds P1605 (迷宫), P1601 (A+B Problem 高精度版) and P1219 (八皇后)
#     as three new rows (32, 33, 34) at the bottom of the problem list.
# ---------------------------------------------------------------------------

# Seed the new rows with the formatting already used by the row above
# (row 31) so that number formats / wrap styles line up with the rest of
# the sheet, then overwrite the individual cell values.
$ws.Range("A31:H31").Copy() | Out-Null
$ws.Range("A32:H34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 32 : P1605 / 迷宫 ---------------------------------------------------
$ws.Range("A32").Value = "P1605"
$ws.Range("B32").Value = "迷宫"
$ws.Range("C32").Value = "AC"
$ws.Range("D32").Value = "普及-"
$ws.Range("E32").Value = "DFS"
$ws.Range("F32").Value = "基本的DFS"
$ws.Range("G32").Value = 43817
$ws.Range("H32").Value = 43817

# --- Row 33 : P1601 / A+B Problem（高精度版） --------------------------------
$ws.Range("A33").Value = "P1601"
$ws.Range("B33").Value = "A+B Problem（高精度版）"
$ws.Range("C33").Value = "AC"
$ws.Range("D33").Value = "普及-"
$ws.Range("E33").Value = "高精度"
$ws.Range("F33").Value = "高精度加法模板"
$ws.Range("G33").Value = 43818
$ws.Range("H33").Value = 43818

# F33 uses the bold / green "hint" style (same one already used on F21)
# rather than the plain wrapped style copied from row 31.
$ws.Range("F21").Copy() | Out-Null
$ws.Range("F33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F33").Value = "高精度加法模板"

# --- Row 34 : P1219 / 八皇后 --------------------------------------------------
$ws.Range("A34").Value = "P1219"
$ws.Range("B34").Value = "八皇后"
$ws.Range("C34").Value = "AC"
$ws.Range("D34").Value = "普及/提高-"
$ws.Range("E34").Value = "DFS"

$hintPart1 = "逐行处理（每行一层递归），检查该行中有无格未被占用，然后把该格序号填入，退出递归时撤销序号`n"
$hintPart2 = "（PS：最后一个点n=13谜之超时，只能打表过？）"
$ws.Range("F34").Value = $hintPart1 + $hintPart2

# Second run of the hint text is rendered a little dimmer (25% darker
# "Background 1" = RGB BFBFBF) to read as a secondary remark.
$dimGray = 12566463
$run2 = $ws.Range("F34").Characters($hintPart1.Length + 1, $hintPart2.Length)
$run2.Font.Color = $dimGray
$run2.Font.Name = "等线"

$ws.Range("G34").Value = 43807
$ws.Range("H34").Value = 43829

# Row 34 wraps onto three lines once the long hint is in place.
$ws.Rows(34).RowHeight = 41.4

# --- View state: scroll down and leave the selection near the new rows ------
$ws.Range("F38:F39").Select() | Out-Null
try { $excel.ActiveWindow.ScrollRow = 26 } catch { }
try { $excel.ActiveWindow.ScrollColumn = 1 } catch { }
